# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 705
    3  = 37
    4  = 537
    5  = 41
    9  = 4466
    10 = 4353
    11 = 7
    12 = 14
    13 = 135
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
